# Daily IST report: add 2026-02-18 column to the submissions matrix.
# Inserts a new date column before the existing "total_files" column,
# fills it with that day's per-person submission counts, and refreshes
# the rolling "total_files" / "unique_days" summary columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 109
$newDateCol = 7      # G: new "2026-02-18" column (old total_files column shifts to H)
$totalCol   = 8      # H: total_files (after shift)
$uniqueCol  = 9      # I: unique_days (after shift)

# 1) Insert a new column at G; everything from G onward (total_files,
#    unique_days) shifts one column to the right.
$ws.Columns("G:G").Insert()

# 2) New column should be the same width as the old "total_files" column
#    used to be (12 chars), matching the updated <cols> block.
$ws.Columns("G:G").ColumnWidth = 11.17

# 3) Header: G1 gets the new date, styled like the other date headers
#    (bold + centered, no fill) rather than the totals header style.
#    Force text so Excel doesn't auto-convert the date-shaped string to
#    a date serial, then restore the date-header look from F1.
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "2026-02-18"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 4) Per-person file counts submitted on 2026-02-18.
$dayCounts = @{2=1; 3=1; 4=0; 5=1; 6=0; 7=1; 8=1; 9=1; 10=1; 11=0; 12=1; 13=1; 14=1; 15=0; 16=0; 17=1; 18=1; 19=1; 20=1; 21=0; 22=0; 23=0; 24=1; 25=1; 26=0; 27=3; 28=0; 29=0; 30=0; 31=1; 32=1; 33=1; 34=1; 35=1; 36=0; 37=1; 38=0; 39=0; 40=1; 41=1; 42=1; 43=21; 44=0; 45=1; 46=1; 47=1; 48=0; 49=1; 50=0; 51=0; 52=1; 53=0; 54=0; 55=1; 56=0; 57=0; 58=1; 59=0; 60=1; 61=1; 62=0; 63=0; 64=1; 65=0; 66=0; 67=0; 68=0; 69=0; 70=1; 71=0; 72=0; 73=0; 74=13; 75=0; 76=1; 77=0; 78=0; 79=0; 80=0; 81=0; 82=1; 83=0; 84=1; 85=0; 86=1; 87=0; 88=0; 89=0; 90=0; 91=0; 92=0; 93=0; 94=0; 95=30; 96=0; 97=0; 98=0; 99=0; 100=1; 101=1; 102=0; 103=1; 104=12; 105=1; 106=1; 107=0; 108=1; 109=0}

for ($r = 2; $r -le $lastRow; $r++) {
    $g = $dayCounts[$r]
    $ws.Cells.Item($r, $newDateCol).Value = $g

    # Recompute the rolling totals across all date columns (D..G).
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2

    $total = $d + $e + $f + $g
    $ws.Cells.Item($r, $totalCol).Value = $total

    $uniqueDays = 0
    if ($d -gt 0) { $uniqueDays++ }
    if ($e -gt 0) { $uniqueDays++ }
    if ($f -gt 0) { $uniqueDays++ }
    if ($g -gt 0) { $uniqueDays++ }
    $ws.Cells.Item($r, $uniqueCol).Value = $uniqueDays
}
